$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update two cell comments on the "Organ" sheet
# ---------------------------------------------------------------------------
$organ = $wb.Worksheets.Item("Organ")

$organ.Range("M1").Comment.Text("The unit of measurement of weight")

$v1Text = "(Required) The string that serves as the definitive identifier for the metadata`nschema version and is readily interpretable by computers for data validation and`nprocessing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9"
$organ.Range("V1").Comment.Text($v1Text)

# ---------------------------------------------------------------------------
# 2. Shrink the warm_ischemic_time_unit / cold_ischemic_time_unit lookup
#    lists from {hour, month, year, day, minute} down to {hour, minute}
# ---------------------------------------------------------------------------
foreach ($name in @("warm_ischemic_time_unit", "cold_ischemic_time_unit")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Item(2, 1).Value = "minute"
    $ws.Cells.Item(2, 2).Value = "http://purl.obolibrary.org/obo/UO_0000031"
    $ws.Range("A3:B5").ClearContents()
}

# ---------------------------------------------------------------------------
# 3. Shrink the weight_unit lookup list from {ng, ug, mg, kg, g} down to
#    {kg, g}
# ---------------------------------------------------------------------------
$weightUnit = $wb.Worksheets.Item("weight_unit")
$weightUnit.Cells.Item(1, 1).Value = "kg"
$weightUnit.Cells.Item(1, 2).Value = "http://purl.obolibrary.org/obo/UO_0000009"
$weightUnit.Cells.Item(2, 1).Value = "g"
$weightUnit.Cells.Item(2, 2).Value = "http://purl.obolibrary.org/obo/UO_0000021"
$weightUnit.Range("A3:B5").ClearContents()

# ---------------------------------------------------------------------------
# 4. Shrink the height_unit / width_unit / length_unit lookup lists from
#    {mm, um, cm, nm} down to {mm, cm}
# ---------------------------------------------------------------------------
foreach ($name in @("height_unit", "width_unit", "length_unit")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Item(2, 1).Value = "cm"
    $ws.Cells.Item(2, 2).Value = "http://purl.obolibrary.org/obo/UO_0000015"
    $ws.Range("A3:B4").ClearContents()
}

# ---------------------------------------------------------------------------
# 5. Reorder/grow the volume_unit lookup list from {cm^3, mm^3, um^3} to
#    {cm^3, um^3, mm^3, ml}
# ---------------------------------------------------------------------------
$volumeUnit = $wb.Worksheets.Item("volume_unit")
$volumeUnit.Cells.Item(2, 1).Value = "um^3"
$volumeUnit.Cells.Item(2, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000112"
$volumeUnit.Cells.Item(3, 1).Value = "mm^3"
$volumeUnit.Cells.Item(3, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000136"
$volumeUnit.Cells.Item(4, 1).Value = "ml"
$volumeUnit.Cells.Item(4, 2).Value = "http://purl.obolibrary.org/obo/UO_0000098"

# ---------------------------------------------------------------------------
# 6. Update the data-validation list ranges on the "Organ" sheet so each
#    dropdown points at the resized lookup ranges
# ---------------------------------------------------------------------------
$organ.Range("H2:H1001").Validation.Delete()
$organ.Range("H2:H1001").Validation.Add(3, 1, 1, "'warm_ischemic_time_unit'!$A$1:$A$2")

$organ.Range("J2:J1001").Validation.Delete()
$organ.Range("J2:J1001").Validation.Add(3, 1, 1, "'cold_ischemic_time_unit'!$A$1:$A$2")

$organ.Range("M2:M1001").Validation.Delete()
$organ.Range("M2:M1001").Validation.Add(3, 1, 1, "'weight_unit'!$A$1:$A$2")

$organ.Range("O2:O1001").Validation.Delete()
$organ.Range("O2:O1001").Validation.Add(3, 1, 1, "'height_unit'!$A$1:$A$2")

$organ.Range("Q2:Q1001").Validation.Delete()
$organ.Range("Q2:Q1001").Validation.Add(3, 1, 1, "'width_unit'!$A$1:$A$2")

$organ.Range("S2:S1001").Validation.Delete()
$organ.Range("S2:S1001").Validation.Add(3, 1, 1, "'length_unit'!$A$1:$A$2")

$organ.Range("U2:U1001").Validation.Delete()
$organ.Range("U2:U1001").Validation.Add(3, 1, 1, "'volume_unit'!$A$1:$A$4")

# ---------------------------------------------------------------------------
# 7. Bump the pav:createdOn timestamp on the .metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Cells.Item(2, 3).Value = "2023-09-08T20:50:05-07:00"
